# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 805 (pushing the existing rows 805-882
# down to 806-883) and populate it with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 805..882 down to 806..883, leaving a blank row 805 in place.
$ws.Rows.Item(805).Insert()

# Populate the newly inserted row 805 with the new record.
$ws.Range("A805").Value = 8
$ws.Range("B805").Value = "Terminal La Palmera de La Serena"
$ws.Range("C805").Value = "Coquimbo"
$ws.Range("D805").Value = 45132
$ws.Range("E805").Value = 4
$ws.Range("F805").Value = 100112024
$ws.Range("G805").Value = "Choclo"
$ws.Range("H805").Value = "Dulce o Americano"
$ws.Range("I805").Value = "Primera"
$ws.Range("J805").Value = 400
$ws.Range("K805").Value = 33000
$ws.Range("L805").Value = 34000
$ws.Range("M805").Value = 33500
$ws.Range("N805").Value = "`$/malla 70 unidades"
$ws.Range("O805").Value = "Región de Arica y Parinacota"
$ws.Range("P805").Value = 479
$ws.Range("Q805").Value = 70
$ws.Range("R805").Value = "Hortaliza"
